# Applies the "6.0.0" release update to the
# StructureDefinition-major-diagnostic-category workbook.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Metadata" (sheet1) ----
$ws1 = $wb.Worksheets.Item(1)

# Version bump
$ws1.Range("B3").Value = "6.0.0"

# Date bump
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value that used to be blank
$ws1.Range("B9").Value = "Alvearie Team"

# The old row 10 ("Contact" / "No display for ContactDetail") becomes
# the new "Jurisdiction" row.
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# The old row 11 was a duplicate "Contact" / "No display for ContactDetail"
# row; remove it entirely so everything below shifts up by one row.
$ws1.Rows.Item(11).Delete()

# ---- Sheet "Elements" (sheet2) ----
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("K2").Value = "Major Diagnostic Category"
$ws2.Range("L2").Value = "Body system or disease related groupings of clinical conditions, based on diagnosis codes"
